# Force-refresh the "Positions" sheet header row with new, more
# descriptive column labels (and the newly-tracked trailing columns),
# then reset the sheet's active selection back to A1.

$wb = $excel.ActiveWorkbook

# Remember whatever sheet is active so we can restore it once we're done
# poking at "Positions" -- we only want to touch that sheet's own
# selection/header state, not steal focus from the workbook's current tab.
$originalActiveSheet = $wb.ActiveSheet.Name

$ws = $wb.Worksheets.Item("Positions")
$ws.Activate()

$headers = @("Ticker", "Broker", "Currency", "Amount", "Cost Basis", "Unit Cost Basis", "Last Price", "Market Value", "Unreal. PnL", "Real. PnL", "Active")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Force the header update to "stick" by resetting the selection to A1.
[void]$ws.Range("A1").Select()

if ($originalActiveSheet -ne $ws.Name) {
    [void]$wb.Worksheets.Item($originalActiveSheet).Activate()
}
